$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.576.93'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -4.50%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.331.05'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -5.28%  '

# Row 4
$ws.Range('E4').Value = '  -0.08%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '549.15'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.19%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '171.14'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -7.50%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.611'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -4.68%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.319.41'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -5.45%  '

# Row 9
$ws.Range('E9').Value = '  +0.02%  '

# Row 10
$ws.Range('E10').Value = '  -4.54%  '

# Row 11
$ws.Range('E11').Value = '  -3.56%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.29'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.88%  '

# Row 13
$ws.Range('E13').Value = '  -5.90%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.88'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.66%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.857.66'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.60%  '

# Row 16
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.325.71'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -5.48%  '

# Row 17
$ws.Range('E17').Value = '  -5.02%  '

# Row 18
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.116'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -4.20%  '

# Row 19
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.63'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.06%  '

# Row 20
$ws.Range('B20').Value = 'WrappedBTC'
$ws.Range('C20').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '63.509.41'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.64%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.968'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.61%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '407.42'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.16%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.05'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.35%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.30'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.62%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.27'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.95%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '82.85'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.45%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.56'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.02%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.72'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -7.88%  '

# Row 29
$ws.Range('E29').Value = '  -6.25%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '29.04'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.71%  '

# Row 31
$ws.Range('E31').Value = '  -7.02%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.29'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.15%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '572.08'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -8.95%  '

# Row 34
$ws.Range('E34').Value = '  -5.75%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '57.23'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.34%  '

# Row 36
$ws.Range('E36').Value = '  +0.08%  '

# Row 37
$ws.Range('E37').Value = '  -2.23%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '35.20'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -8.06%  '

# Row 39
$ws.Range('E39').Value = '  +0.70%  '

# Row 40
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.164.48'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.90%  '

# Row 41
$ws.Range('B41').Value = 'PEPE'
$ws.Range('C41').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0₃0735'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -9.93%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.366'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.71%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.999'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.22%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.82'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.02%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.20'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.82%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.44'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -7.46%  '

# Row 47
$ws.Range('E47').Value = '  -4.60%  '

# Row 48
$ws.Range('E48').Value = '  -4.80%  '

# Row 49
$ws.Range('E49').Value = '  -4.43%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.59'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.42%  '

# Row 51
$ws.Range('E51').Value = '  -6.47%  '
